# ContosoLearn Competitor SWOT: translate body text from English to Spanish.
#
# Strategy: locate each paragraph's current range fresh before every edit
# (offsets shift as text is replaced), replace the bold "header" word
# (Strengths/Weaknesses/Opportunities/Threats) in place, then replace the
# descriptive sentence(s) that follow. Where the source was edited with a
# translation tool that emitted each sentence as its own run, we reproduce
# that by nudging the font of the newly-inserted sentence (same visual
# font, Aptos) which forces the engine to split a fresh <w:r> at that
# boundary, and we tag the new Spanish text as es-ES.

$d = $word.ActiveDocument

function Replace-Text($doc, $startAnchor, $endAnchor, $oldText, $newText) {
    # Locate-then-assign (rather than Find.Execute's built-in Replace)
    # so that AutoCorrect "smart quotes" never mangles literal straight
    # quote characters in the replacement text.
    $r = $doc.Range($startAnchor, $endAnchor)
    $ok = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find failed for: $oldText"
    }
    $hit = $doc.Range($r.Start, $r.End)
    $hit.Text = $newText
}

function Mark-Run($doc, $startAnchor, $endAnchor, $text) {
    # Locate $text within [$startAnchor, $endAnchor) and apply a direct
    # (no-op-value) font change so the engine splits a new run exactly at
    # this boundary, then tag the run as Spanish.
    $r = $doc.Range($startAnchor, $endAnchor)
    $ok = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find failed for: $text"
    }
    $hit = $doc.Range($r.Start, $r.End)
    $hit.Font.Name = "Aptos"
    $hit.LanguageID = "es-ES"
}

# ---------------------------------------------------------------------
# Paragraph 1 (title)
# ---------------------------------------------------------------------
$p = $d.Paragraphs(1)
Replace-Text $d $p.Range.Start $p.Range.End "ContosoLearn Competitor SWOT" "DAFO de competidores ContosoLearn"

# ---------------------------------------------------------------------
# Paragraph 2: "Fabrikam Learning:" header (unchanged text)
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Paragraph 3: Strengths (Fabrikam)
# ---------------------------------------------------------------------
$p = $d.Paragraphs(3)
$pStart = $p.Range.Start
Replace-Text $d $pStart $p.Range.End "Strengths:" "Fortalezas:"

$pStart = $d.Paragraphs(3).Range.Start
Replace-Text $d $pStart $d.Paragraphs(3).Range.End `
    " Fabrikam Learning provides a comprehensive set of analytics and reporting tools." `
    " Fabrikam Learning proporciona un conjunto completo de herramientas de análisis e informes."

$pStart = $d.Paragraphs(3).Range.Start
Replace-Text $d $pStart $d.Paragraphs(3).Range.End `
    " It ensures the continuous monitoring of teaching and learning activities, as well as pinpointing problematic areas that need to be addressed." `
    " Garantiza la supervisión continua de las actividades de enseñanza y aprendizaje, así como la identificación de áreas problemáticas que deben abordarse."

$pStart = $d.Paragraphs(3).Range.Start
$pEnd = $d.Paragraphs(3).Range.End
Mark-Run $d $pStart $pEnd "Fabrikam Learning proporciona un conjunto completo de herramientas de análisis e informes."
$pEnd = $d.Paragraphs(3).Range.End
Mark-Run $d $pStart $pEnd " "
$pEnd = $d.Paragraphs(3).Range.End
Mark-Run $d $pStart $pEnd "Garantiza la supervisión continua de las actividades de enseñanza y aprendizaje, así como la identificación de áreas problemáticas que deben abordarse."

# ---------------------------------------------------------------------
# Paragraph 4: Weaknesses (Fabrikam)
# ---------------------------------------------------------------------
$p = $d.Paragraphs(4)
$pStart = $p.Range.Start
Replace-Text $d $pStart $p.Range.End "Weaknesses:" "Puntos débiles:"

$pStart = $d.Paragraphs(4).Range.Start
Replace-Text $d $pStart $d.Paragraphs(4).Range.End `
    " While Fabrikam Learning has robust reporting capabilities, it might be overwhelming for some users due to its comprehensive nature." `
    " aunque Fabrikam Learning tiene funcionalidades de informes sólidas, puede ser abrumador para algunos usuarios debido a su naturaleza completa."

$pStart = $d.Paragraphs(4).Range.Start
$pEnd = $d.Paragraphs(4).Range.End
Mark-Run $d $pStart $pEnd "aunque Fabrikam Learning tiene funcionalidades de informes sólidas, puede ser abrumador para algunos usuarios debido a su naturaleza completa."

# ---------------------------------------------------------------------
# Paragraph 5: Opportunities (Fabrikam)
# ---------------------------------------------------------------------
$p = $d.Paragraphs(5)
$pStart = $p.Range.Start
Replace-Text $d $pStart $p.Range.End "Opportunities:" "Oportunidades:"

$pStart = $d.Paragraphs(5).Range.Start
Replace-Text $d $pStart $d.Paragraphs(5).Range.End `
    " There is a growing demand for personalized learning experiences and data-driven recommendations." `
    " hay una creciente demanda de experiencias de aprendizaje personalizadas y recomendaciones controladas por datos."

$pStart = $d.Paragraphs(5).Range.Start
Replace-Text $d $pStart $d.Paragraphs(5).Range.End `
    " Fabrikam Learning can leverage its robust analytics and reporting tools to meet this demand." `
    " Fabrikam Learning puede aprovechar sus sólidas herramientas de análisis e informes para satisfacer esta demanda."

$pStart = $d.Paragraphs(5).Range.Start
$pEnd = $d.Paragraphs(5).Range.End
Mark-Run $d $pStart $pEnd "hay una creciente demanda de experiencias de aprendizaje personalizadas y recomendaciones controladas por datos."
$pEnd = $d.Paragraphs(5).Range.End
Mark-Run $d $pStart $pEnd " "
$pEnd = $d.Paragraphs(5).Range.End
Mark-Run $d $pStart $pEnd "Fabrikam Learning puede aprovechar sus sólidas herramientas de análisis e informes para satisfacer esta demanda."

# ---------------------------------------------------------------------
# Paragraph 6: Threats (Fabrikam)
# ---------------------------------------------------------------------
$p = $d.Paragraphs(6)
$pStart = $p.Range.Start
Replace-Text $d $pStart $p.Range.End "Threats:" "Amenazas:"

$pStart = $d.Paragraphs(6).Range.Start
Replace-Text $d $pStart $d.Paragraphs(6).Range.End `
    " The eLearning market is highly competitive with many players offering similar features." `
    " el mercado de eLearning es altamente competitivo con muchos jugadores que ofrecen características similares."

$pStart = $d.Paragraphs(6).Range.Start
Replace-Text $d $pStart $d.Paragraphs(6).Range.End `
    " Fabrikam Learning needs to continuously innovate to stay ahead." `
    " Fabrikam Learning debe innovar continuamente para mantenerse a la vanguardia."

$pStart = $d.Paragraphs(6).Range.Start
$pEnd = $d.Paragraphs(6).Range.End
Mark-Run $d $pStart $pEnd "el mercado de eLearning es altamente competitivo con muchos jugadores que ofrecen características similares."
$pEnd = $d.Paragraphs(6).Range.End
Mark-Run $d $pStart $pEnd " "
$pEnd = $d.Paragraphs(6).Range.End
Mark-Run $d $pStart $pEnd "Fabrikam Learning debe innovar continuamente para mantenerse a la vanguardia."

# ---------------------------------------------------------------------
# Paragraph 7: "AdatumLearn:" header (unchanged text)
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Paragraph 8: Strengths (AdatumLearn)
# ---------------------------------------------------------------------
$p = $d.Paragraphs(8)
$pStart = $p.Range.Start
Replace-Text $d $pStart $p.Range.End "Strengths:" "Fortalezas:"

$pStart = $d.Paragraphs(8).Range.Start
Replace-Text $d $pStart $d.Paragraphs(8).Range.End `
    " AdatumLearn offers courses on business analysis techniques such as MOST and SWOT." `
    " AdatumLearn ofrece cursos sobre técnicas de análisis de negocios como MOST y DAFO."

$pStart = $d.Paragraphs(8).Range.Start
Replace-Text $d $pStart $d.Paragraphs(8).Range.End `
    " This shows their commitment to providing valuable content to their users." `
    " Esto muestra su compromiso de proporcionar contenido valioso a sus usuarios."

$pStart = $d.Paragraphs(8).Range.Start
$pEnd = $d.Paragraphs(8).Range.End
Mark-Run $d $pStart $pEnd "AdatumLearn ofrece cursos sobre técnicas de análisis de negocios como MOST y DAFO."
$pEnd = $d.Paragraphs(8).Range.End
Mark-Run $d $pStart $pEnd " "
$pEnd = $d.Paragraphs(8).Range.End
Mark-Run $d $pStart $pEnd "Esto muestra su compromiso de proporcionar contenido valioso a sus usuarios."

# ---------------------------------------------------------------------
# Paragraph 9: Weaknesses (AdatumLearn)
# ---------------------------------------------------------------------
$p = $d.Paragraphs(9)
$pStart = $p.Range.Start
Replace-Text $d $pStart $p.Range.End "Weaknesses:" "Puntos débiles:"

$pStart = $d.Paragraphs(9).Range.Start
Replace-Text $d $pStart $d.Paragraphs(9).Range.End `
    " The information provided in their courses is a compilation of third-party generated information." `
    " la información proporcionada en sus cursos es una compilación de información generada por terceros."

$pStart = $d.Paragraphs(9).Range.Start
Replace-Text $d $pStart $d.Paragraphs(9).Range.End `
    " This might not be as valuable as original content." `
    " Esto podría no ser tan valioso como el contenido original."

$pStart = $d.Paragraphs(9).Range.Start
$pEnd = $d.Paragraphs(9).Range.End
Mark-Run $d $pStart $pEnd "la información proporcionada en sus cursos es una compilación de información generada por terceros."
$pEnd = $d.Paragraphs(9).Range.End
Mark-Run $d $pStart $pEnd " "
$pEnd = $d.Paragraphs(9).Range.End
Mark-Run $d $pStart $pEnd "Esto podría no ser tan valioso como el contenido original."

# ---------------------------------------------------------------------
# Paragraph 10: Opportunities (AdatumLearn)
# ---------------------------------------------------------------------
$p = $d.Paragraphs(10)
$pStart = $p.Range.Start
Replace-Text $d $pStart $p.Range.End "Opportunities:" "Oportunidades:"

$pStart = $d.Paragraphs(10).Range.Start
Replace-Text $d $pStart $d.Paragraphs(10).Range.End `
    " AdatumLearn can create more original content to provide unique value to their users." `
    " AdatumLearn puede crear contenido más original para proporcionar un valor único a sus usuarios."

$pStart = $d.Paragraphs(10).Range.Start
Replace-Text $d $pStart $d.Paragraphs(10).Range.End `
    " They can also expand their course offerings to cover more topics." `
    " También puede ampliar sus ofertas de cursos para tratar más temas."

$pStart = $d.Paragraphs(10).Range.Start
$pEnd = $d.Paragraphs(10).Range.End
Mark-Run $d $pStart $pEnd "AdatumLearn puede crear contenido más original para proporcionar un valor único a sus usuarios."
$pEnd = $d.Paragraphs(10).Range.End
Mark-Run $d $pStart $pEnd " "
$pEnd = $d.Paragraphs(10).Range.End
Mark-Run $d $pStart $pEnd "También puede ampliar sus ofertas de cursos para tratar más temas."

# ---------------------------------------------------------------------
# Paragraph 11: Threats (AdatumLearn)
# ---------------------------------------------------------------------
$p = $d.Paragraphs(11)
$pStart = $p.Range.Start
Replace-Text $d $pStart $p.Range.End "Threats:" "Amenazas:"

$pStart = $d.Paragraphs(11).Range.Start
Replace-Text $d $pStart $d.Paragraphs(11).Range.End `
    " Like Fabrikam Learning, AdatumLearn also faces stiff competition in the eLearning market." `
    " al igual que Fabrikam Learning, AdatumLearn también se enfrenta a una competencia rígida en el mercado de eLearning."

$pStart = $d.Paragraphs(11).Range.Start
Replace-Text $d $pStart $d.Paragraphs(11).Range.End `
    " They need to continuously improve their offerings to stay competitive.`"" `
    " Necesita mejorar continuamente su oferta para mantenerse competitivo`"."

$pStart = $d.Paragraphs(11).Range.Start
$pEnd = $d.Paragraphs(11).Range.End
Mark-Run $d $pStart $pEnd "al igual que Fabrikam Learning, AdatumLearn también se enfrenta a una competencia rígida en el mercado de eLearning."
$pEnd = $d.Paragraphs(11).Range.End
Mark-Run $d $pStart $pEnd " "
$pEnd = $d.Paragraphs(11).Range.End
Mark-Run $d $pStart $pEnd "Necesita mejorar continuamente su oferta para mantenerse competitivo`"."

Write-Host "Done."
